$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -92.3178
$ws.Range("B2").Value = -92.2627

$ws.Range("A3").Value = 46.5152
$ws.Range("B3").Value = 46.5536

$ws.Range("A4").Value = -91.9127
$ws.Range("B4").Value = -91.9682

$ws.Range("A5").Value = 46.7961
$ws.Range("B5").Value = 46.7578
